# Edit script: "added slide on standards"
#
# Summary of changes applied to the deck:
#   1. Cosmetic run-merges (text content unchanged, just collapsed
#      into single runs) on a few existing slides.
#   2. The "Sharing and re-using" slide is moved earlier (from position 5
#      to position 4, ahead of "Versioning").
#   3. A brand-new "Standardizing" slide is inserted right after it
#      (new position 5), pushing "Versioning" and everything after it
#      one slot later.

$p = $ppt.ActivePresentation

function Set-ParaText($textRange, $paraIndex, $newText) {
    # Force a full run-rebuild (collapsing any pre-existing run-splits)
    # by first writing a sentinel value, then the real text. A direct
    # write of text that is already equal to the paragraph's current
    # (run-concatenated) text is treated as a no-op by the host and
    # would NOT collapse multiple runs into one.
    $para = $textRange.Paragraphs($paraIndex)
    $para.Text = "zzz__tmp__zzz"
    $textRange.Paragraphs($paraIndex).Text = $newText
}

# --- 1. Cosmetic run-merge fixes (text content is unchanged) ---------

# Slide 3 "Normalizing and integrating"
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(2).TextFrame.TextRange
Set-ParaText $tr3 1 "Your data are linked to things in the real world"
Set-ParaText $tr3 5 "Links to trees and studies (TreeBASE)"

# Slide 5 "Sharing and re-using" (still at position 5 at this point)
$s5 = $p.Slides.Item(5)
$tr5 = $s5.Shapes.Item(2).TextFrame.TextRange
Set-ParaText $tr5 2 "In a “Big Science”, we can’t do all on our own"

# Slide 10 "Stay in touch!"
$s10 = $p.Slides.Item(10)
$tr10 = $s10.Shapes.Item(3).TextFrame.TextRange
Set-ParaText $tr10 3 "@rvosa, @nescent, @treebase, @phylofoundation, @museumnaturalis, etc."

# --- 2. Reorder: "Sharing and re-using" moves ahead of "Versioning" ---

$p.Slides.Item(5).MoveTo(4)

# --- 3. Insert new "Standardizing" slide at position 5 ----------------

$newSlide = $p.Slides.Add(5, 2)   # 2 = "Title and Content" custom layout

$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Standardizing"

$body = $newSlide.Shapes.Item(2).TextFrame.TextRange
$body.Text = "Data standards enable sharing, re-using, and integrating`rUnderstand the standards that you work with, and follow them:`rNeXML, phyloXML, NEXUS, NEWICK, FASTA, etc."
$body.Paragraphs(1).IndentLevel = 1
$body.Paragraphs(2).IndentLevel = 1
$body.Paragraphs(3).IndentLevel = 2
